# Fix Training Data Issue (#48)
# Correct team-stat values and normalize the Date column format
# (data had been captured one calendar day off; also reformat
#  "11-20-2008-09" -> "2008-11-20").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AD2").Value = 14
$ws.Range("AF2").Value = 6
$ws.Range("AG2").Value = 9
$ws.Range("AI2").Value = 12
$ws.Range("AO2").Value = 19
$ws.Range("AT2").Value = 14
$ws.Range("AU2").Value = 14
$ws.Range("AV2").Value = 7
$ws.Range("AZ2").Value = 13
$ws.Range("BC2").Value = 13
$ws.Range("BF2").Value = "'2008-11-20"

# Row 3
$ws.Range("I3").Value = 33.8
$ws.Range("J3").Value = 75.5
$ws.Range("K3").Value = 0.447
$ws.Range("L3").Value = 5.1
$ws.Range("N3").Value = 0.311
$ws.Range("O3").Value = 23
$ws.Range("P3").Value = 30.3
$ws.Range("Q3").Value = 0.758
$ws.Range("R3").Value = 10.8
$ws.Range("S3").Value = 32.6
$ws.Range("T3").Value = 43.4
$ws.Range("U3").Value = 20.3
$ws.Range("V3").Value = 17.3
$ws.Range("W3").Value = 8.6
$ws.Range("X3").Value = 5.6
$ws.Range("Y3").Value = 3.9
$ws.Range("Z3").Value = 24.1
$ws.Range("AA3").Value = 25.3
$ws.Range("AB3").Value = 95.6
$ws.Range("AC3").Value = 4.9
$ws.Range("AH3").Value = 10
$ws.Range("AK3").Value = 12
$ws.Range("AL3").Value = 20
$ws.Range("AN3").Value = 24
$ws.Range("AQ3").Value = 18
$ws.Range("AS3").Value = 5
$ws.Range("AT3").Value = 10
$ws.Range("AU3").Value = 13
$ws.Range("AX3").Value = 9
$ws.Range("AY3").Value = 4
$ws.Range("AZ3").Value = 28
$ws.Range("BB3").Value = 22
$ws.Range("BF3").Value = "'2008-11-20"

# Row 4
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 26
$ws.Range("AN4").Value = 23
$ws.Range("AO4").Value = 10
$ws.Range("AQ4").Value = 11
$ws.Range("AV4").Value = 19
$ws.Range("AY4").Value = 30
$ws.Range("BA4").Value = 13
$ws.Range("BF4").Value = "'2008-11-20"

# Row 5
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 5
$ws.Range("G5").Value = 0.417
$ws.Range("J5").Value = 84.2
$ws.Range("K5").Value = 0.415
$ws.Range("M5").Value = 14.8
$ws.Range("N5").Value = 0.339
$ws.Range("O5").Value = 21.8
$ws.Range("P5").Value = 27.6
$ws.Range("Q5").Value = 0.792
$ws.Range("R5").Value = 13.3
$ws.Range("S5").Value = 30.8
$ws.Range("T5").Value = 44
$ws.Range("U5").Value = 18.5
$ws.Range("V5").Value = 15.3
$ws.Range("X5").Value = 5.6
$ws.Range("Y5").Value = 6.8
$ws.Range("Z5").Value = 22.8
$ws.Range("AA5").Value = 22.3
$ws.Range("AB5").Value = 96.7
$ws.Range("AC5").Value = -3
$ws.Range("AD5").Value = 3
$ws.Range("AE5").Value = 15
$ws.Range("AG5").Value = 20
$ws.Range("AJ5").Value = 8
$ws.Range("AM5").Value = 24
$ws.Range("AO5").Value = 4
$ws.Range("AP5").Value = 8
$ws.Range("AQ5").Value = 5
$ws.Range("AT5").Value = 7
$ws.Range("AU5").Value = 23
$ws.Range("AV5").Value = 21
$ws.Range("AX5").Value = 9
$ws.Range("AY5").Value = 29
$ws.Range("AZ5").Value = 22
$ws.Range("BA5").Value = 11
$ws.Range("BB5").Value = 18
$ws.Range("BC5").Value = 22
$ws.Range("BF5").Value = "'2008-11-20"

# Row 6
$ws.Range("D6").Value = 12
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 0.75
$ws.Range("I6").Value = 36
$ws.Range("J6").Value = 76.3
$ws.Range("K6").Value = 0.472
$ws.Range("L6").Value = 6.7
$ws.Range("M6").Value = 19
$ws.Range("N6").Value = 0.351
$ws.Range("O6").Value = 21.8
$ws.Range("P6").Value = 28.1
$ws.Range("Q6").Value = 0.777
$ws.Range("R6").Value = 10.7
$ws.Range("S6").Value = 29.8
$ws.Range("U6").Value = 20
$ws.Range("V6").Value = 13.7
$ws.Range("X6").Value = 5.9
$ws.Range("Y6").Value = 3.3
$ws.Range("Z6").Value = 21.8
$ws.Range("AA6").Value = 22.8
$ws.Range("AB6").Value = 100.5
$ws.Range("AC6").Value = 7.7
$ws.Range("AD6").Value = 3
$ws.Range("AF6").Value = 3
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 26
$ws.Range("AL6").Value = 10
$ws.Range("AM6").Value = 9
$ws.Range("AN6").Value = 15
$ws.Range("AQ6").Value = 10
$ws.Range("AS6").Value = 20
$ws.Range("AU6").Value = 17
$ws.Range("AV6").Value = 9
$ws.Range("AX6").Value = 7
$ws.Range("AY6").Value = 1
$ws.Range("AZ6").Value = 14
$ws.Range("BA6").Value = 7
$ws.Range("BF6").Value = "'2008-11-20"

# Row 7
$ws.Range("AH7").Value = 10
$ws.Range("AJ7").Value = 5
$ws.Range("AO7").Value = 18
$ws.Range("AP7").Value = 18
$ws.Range("AQ7").Value = 9
$ws.Range("AT7").Value = 3
$ws.Range("AV7").Value = 10
$ws.Range("AX7").Value = 16
$ws.Range("BF7").Value = "'2008-11-20"

# Row 8
$ws.Range("AF8").Value = 6
$ws.Range("AG8").Value = 6
$ws.Range("AH8").Value = 10
$ws.Range("AL8").Value = 20
$ws.Range("AN8").Value = 22
$ws.Range("AT8").Value = 15
$ws.Range("AV8").Value = 27
$ws.Range("AZ8").Value = 26
$ws.Range("BB8").Value = 9
$ws.Range("BF8").Value = "'2008-11-20"

# Row 9
$ws.Range("D9").Value = 11
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 0.727
$ws.Range("I9").Value = 35.8
$ws.Range("J9").Value = 79.9
$ws.Range("K9").Value = 0.448
$ws.Range("L9").Value = 5.9
$ws.Range("M9").Value = 15.2
$ws.Range("O9").Value = 21.1
$ws.Range("P9").Value = 26.5
$ws.Range("Q9").Value = 0.795
$ws.Range("R9").Value = 12.2
$ws.Range("S9").Value = 30.1
$ws.Range("T9").Value = 42.3
$ws.Range("U9").Value = 20.5
$ws.Range("V9").Value = 13.2
$ws.Range("X9").Value = 5.3
$ws.Range("Z9").Value = 22.5
$ws.Range("AA9").Value = 22.7
$ws.Range("AB9").Value = 98.6
$ws.Range("AC9").Value = 3
$ws.Range("AD9").Value = 14
$ws.Range("AF9").Value = 3
$ws.Range("AG9").Value = 4
$ws.Range("AI9").Value = 14
$ws.Range("AJ9").Value = 14
$ws.Range("AK9").Value = 11
$ws.Range("AP9").Value = 12
$ws.Range("AQ9").Value = 4
$ws.Range("AS9").Value = 17
$ws.Range("AT9").Value = 13
$ws.Range("AU9").Value = 9
$ws.Range("AV9").Value = 6
$ws.Range("AZ9").Value = 19
$ws.Range("BB9").Value = 10
$ws.Range("BC9").Value = 8
$ws.Range("BF9").Value = "'2008-11-20"

# Row 10
$ws.Range("AD10").Value = 14
$ws.Range("AJ10").Value = 2
$ws.Range("AL10").Value = 18
$ws.Range("AS10").Value = 15
$ws.Range("AZ10").Value = 23
$ws.Range("BB10").Value = 2
$ws.Range("BF10").Value = "'2008-11-20"

# Row 11
$ws.Range("AG11").Value = 10
$ws.Range("AH11").Value = 10
$ws.Range("AM11").Value = 14
$ws.Range("AR11").Value = 20
$ws.Range("AV11").Value = 8
$ws.Range("BC11").Value = 12
$ws.Range("BF11").Value = "'2008-11-20"

# Row 12
$ws.Range("AJ12").Value = 10
$ws.Range("AK12").Value = 6
$ws.Range("AN12").Value = 14
$ws.Range("AQ12").Value = 13
$ws.Range("AS12").Value = 4
$ws.Range("AU12").Value = 5
$ws.Range("AV12").Value = 26
$ws.Range("BC12").Value = 7
$ws.Range("BF12").Value = "'2008-11-20"

# Row 13
$ws.Range("AD13").Value = 14
$ws.Range("AI13").Value = 19
$ws.Range("AZ13").Value = 18
$ws.Range("BF13").Value = "'2008-11-20"

# Row 14
$ws.Range("I14").Value = 38.4
$ws.Range("J14").Value = 86.2
$ws.Range("L14").Value = 6.3
$ws.Range("M14").Value = 16.4
$ws.Range("N14").Value = 0.385
$ws.Range("O14").Value = 21.7
$ws.Range("P14").Value = 28.6
$ws.Range("Q14").Value = 0.759
$ws.Range("S14").Value = 33.9
$ws.Range("T14").Value = 47.9
$ws.Range("U14").Value = 21.3
$ws.Range("V14").Value = 14.2
$ws.Range("W14").Value = 10.7
$ws.Range("X14").Value = 6.6
$ws.Range("Y14").Value = 4.7
$ws.Range("Z14").Value = 21
$ws.Range("AA14").Value = 22.8
$ws.Range("AB14").Value = 104.9
$ws.Range("AC14").Value = 13.6
$ws.Range("AI14").Value = 2
$ws.Range("AJ14").Value = 3
$ws.Range("AM14").Value = 15
$ws.Range("AN14").Value = 5
$ws.Range("AP14").Value = 6
$ws.Range("AQ14").Value = 17
$ws.Range("AS14").Value = 3
$ws.Range("AT14").Value = 2
$ws.Range("AU14").Value = 6
$ws.Range("AV14").Value = 13
$ws.Range("AY14").Value = 12
$ws.Range("AZ14").Value = 11
$ws.Range("BA14").Value = 9
$ws.Range("BB14").Value = 3
$ws.Range("BF14").Value = "'2008-11-20"

# Row 15
$ws.Range("AD15").Value = 14
$ws.Range("AE15").Value = 24
$ws.Range("AG15").Value = 25
$ws.Range("AJ15").Value = 19
$ws.Range("AM15").Value = 23
$ws.Range("AN15").Value = 20
$ws.Range("AP15").Value = 14
$ws.Range("AQ15").Value = 13
$ws.Range("AZ15").Value = 23
$ws.Range("BA15").Value = 14
$ws.Range("BC15").Value = 23
$ws.Range("BF15").Value = "'2008-11-20"

# Row 16
$ws.Range("AI16").Value = 11
$ws.Range("AK16").Value = 9
$ws.Range("AP16").Value = 17
$ws.Range("AU16").Value = 11
$ws.Range("AX16").Value = 8
$ws.Range("BC16").Value = 6
$ws.Range("BF16").Value = "'2008-11-20"

# Row 17
$ws.Range("D17").Value = 13
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 0.385
$ws.Range("H17").Value = 49.2
$ws.Range("I17").Value = 35.2
$ws.Range("J17").Value = 81.9
$ws.Range("K17").Value = 0.43
$ws.Range("N17").Value = 0.375
$ws.Range("O17").Value = 19.8
$ws.Range("P17").Value = 26.4
$ws.Range("R17").Value = 13.6
$ws.Range("S17").Value = 29.9
$ws.Range("T17").Value = 43.5
$ws.Range("V17").Value = 16.5
$ws.Range("W17").Value = 5.6
$ws.Range("Z17").Value = 26.2
$ws.Range("AA17").Value = 24.4
$ws.Range("AB17").Value = 95.6
$ws.Range("AC17").Value = -2.5
$ws.Range("AD17").Value = 1
$ws.Range("AF17").Value = 25
$ws.Range("AG17").Value = 23
$ws.Range("AI17").Value = 17
$ws.Range("AJ17").Value = 11
$ws.Range("AK17").Value = 25
$ws.Range("AL17").Value = 19
$ws.Range("AN17").Value = 10
$ws.Range("AO17").Value = 14
$ws.Range("AP17").Value = 13
$ws.Range("AS17").Value = 18
$ws.Range("AT17").Value = 9
$ws.Range("AU17").Value = 8
$ws.Range("AV17").Value = 28
$ws.Range("BA17").Value = 4
$ws.Range("BF17").Value = "'2008-11-20"

# Row 18
$ws.Range("AF18").Value = 25
$ws.Range("AO18").Value = 20
$ws.Range("AQ18").Value = 3
$ws.Range("AZ18").Value = 27
$ws.Range("BC18").Value = 24
$ws.Range("BF18").Value = "'2008-11-20"

# Row 19
$ws.Range("AE19").Value = 24
$ws.Range("AN19").Value = 6
$ws.Range("AO19").Value = 7
$ws.Range("AP19").Value = 9
$ws.Range("AQ19").Value = 7
$ws.Range("AU19").Value = 24
$ws.Range("AV19").Value = 12
$ws.Range("BA19").Value = 16
$ws.Range("BF19").Value = "'2008-11-20"

# Row 20
$ws.Range("AK20").Value = 7
$ws.Range("AM20").Value = 8
$ws.Range("AN20").Value = 9
$ws.Range("AQ20").Value = 12
$ws.Range("AU20").Value = 12
$ws.Range("AV20").Value = 5
$ws.Range("AY20").Value = 3
$ws.Range("BF20").Value = "'2008-11-20"

# Row 21
$ws.Range("AD21").Value = 14
$ws.Range("AI21").Value = 1
$ws.Range("AN21").Value = 8
$ws.Range("AQ21").Value = 20
$ws.Range("AS21").Value = 11
$ws.Range("AY21").Value = 10
$ws.Range("BB21").Value = 1
$ws.Range("BF21").Value = "'2008-11-20"

# Row 22
$ws.Range("AJ22").Value = 6
$ws.Range("AN22").Value = 13
$ws.Range("AQ22").Value = 19
$ws.Range("AT22").Value = 6
$ws.Range("AZ22").Value = 17
$ws.Range("BF22").Value = "'2008-11-20"

# Row 23
$ws.Range("AD23").Value = 14
$ws.Range("AF23").Value = 3
$ws.Range("AI23").Value = 16
$ws.Range("AJ23").Value = 15
$ws.Range("AL23").Value = 4
$ws.Range("AN23").Value = 21
$ws.Range("AT23").Value = 8
$ws.Range("AV23").Value = 11
$ws.Range("AW23").Value = 12
$ws.Range("AY23").Value = 6
$ws.Range("AZ23").Value = 12
$ws.Range("BA23").Value = 8
$ws.Range("BF23").Value = "'2008-11-20"

# Row 24
$ws.Range("AD24").Value = 14
$ws.Range("AJ24").Value = 7
$ws.Range("AN24").Value = 12
$ws.Range("AU24").Value = 17
$ws.Range("AV24").Value = 25
$ws.Range("BB24").Value = 17
$ws.Range("BC24").Value = 11
$ws.Range("BF24").Value = "'2008-11-20"

# Row 25
$ws.Range("D25").Value = 12
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 0.667
$ws.Range("I25").Value = 36.4
$ws.Range("J25").Value = 73
$ws.Range("K25").Value = 0.499
$ws.Range("M25").Value = 16
$ws.Range("N25").Value = 0.38
$ws.Range("O25").Value = 21.6
$ws.Range("P25").Value = 28.9
$ws.Range("Q25").Value = 0.746
$ws.Range("R25").Value = 8
$ws.Range("S25").Value = 31.3
$ws.Range("T25").Value = 39.3
$ws.Range("U25").Value = 19.9
$ws.Range("V25").Value = 16.8
$ws.Range("W25").Value = 6.3
$ws.Range("X25").Value = 5.1
$ws.Range("Y25").Value = 4.4
$ws.Range("Z25").Value = 20.8
$ws.Range("AA25").Value = 23.1
$ws.Range("AB25").Value = 100.5
$ws.Range("AC25").Value = 2.8
$ws.Range("AD25").Value = 3
$ws.Range("AF25").Value = 6
$ws.Range("AG25").Value = 6
$ws.Range("AH25").Value = 10
$ws.Range("AI25").Value = 10
$ws.Range("AJ25").Value = 30
$ws.Range("AN25").Value = 7
$ws.Range("AO25").Value = 8
$ws.Range("AP25").Value = 5
$ws.Range("AQ25").Value = 22
$ws.Range("AR25").Value = 30
$ws.Range("AS25").Value = 12
$ws.Range("AT25").Value = 25
$ws.Range("AU25").Value = 19
$ws.Range("AY25").Value = 9
$ws.Range("BA25").Value = 6
$ws.Range("BB25").Value = 4
$ws.Range("BC25").Value = 9
$ws.Range("BF25").Value = "'2008-11-20"

# Row 26
$ws.Range("D26").Value = 12
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 0.583
$ws.Range("H26").Value = 48.4
$ws.Range("I26").Value = 35.4
$ws.Range("J26").Value = 78.8
$ws.Range("K26").Value = 0.449
$ws.Range("L26").Value = 8.3
$ws.Range("M26").Value = 19.7
$ws.Range("N26").Value = 0.419
$ws.Range("O26").Value = 18.3
$ws.Range("P26").Value = 23.8
$ws.Range("Q26").Value = 0.768
$ws.Range("R26").Value = 12.3
$ws.Range("S26").Value = 27.7
$ws.Range("T26").Value = 39.9
$ws.Range("U26").Value = 20.5
$ws.Range("V26").Value = 12.4
$ws.Range("W26").Value = 7.3
$ws.Range("X26").Value = 5.4
$ws.Range("Y26").Value = 3.3
$ws.Range("Z26").Value = 21.9
$ws.Range("AA26").Value = 21.3
$ws.Range("AB26").Value = 97.3
$ws.Range("AC26").Value = 1.8
$ws.Range("AD26").Value = 3
$ws.Range("AF26").Value = 10
$ws.Range("AG26").Value = 10
$ws.Range("AH26").Value = 10
$ws.Range("AI26").Value = 15
$ws.Range("AJ26").Value = 17
$ws.Range("AK26").Value = 10
$ws.Range("AM26").Value = 6
$ws.Range("AO26").Value = 21
$ws.Range("AP26").Value = 19
$ws.Range("AQ26").Value = 15
$ws.Range("AT26").Value = 22
$ws.Range("AU26").Value = 10
$ws.Range("AW26").Value = 19
$ws.Range("AX26").Value = 12
$ws.Range("AZ26").Value = 16
$ws.Range("BA26").Value = 17
$ws.Range("BB26").Value = 16
$ws.Range("BC26").Value = 14
$ws.Range("BF26").Value = "'2008-11-20"

# Row 27
$ws.Range("AF27").Value = 25
$ws.Range("AH27").Value = 16
$ws.Range("AJ27").Value = 21
$ws.Range("AQ27").Value = 6
$ws.Range("AU27").Value = 16
$ws.Range("AY27").Value = 13
$ws.Range("AZ27").Value = 25
$ws.Range("BF27").Value = "'2008-11-20"

# Row 28
$ws.Range("AD28").Value = 14
$ws.Range("AI28").Value = 18
$ws.Range("AK28").Value = 8
$ws.Range("AM28").Value = 7
$ws.Range("AN28").Value = 11
$ws.Range("AQ28").Value = 23
$ws.Range("AR28").Value = 29
$ws.Range("AS28").Value = 14
$ws.Range("AU28").Value = 14
$ws.Range("BF28").Value = "'2008-11-20"

# Row 29
$ws.Range("AD29").Value = 14
$ws.Range("AI29").Value = 19
$ws.Range("AJ29").Value = 25
$ws.Range("AK29").Value = 5
$ws.Range("AS29").Value = 15
$ws.Range("AT29").Value = 24
$ws.Range("AX29").Value = 11
$ws.Range("BA29").Value = 14
$ws.Range("BF29").Value = "'2008-11-20"

# Row 30
$ws.Range("AF30").Value = 6
$ws.Range("AG30").Value = 6
$ws.Range("AJ30").Value = 20
$ws.Range("AO30").Value = 13
$ws.Range("AP30").Value = 15
$ws.Range("AQ30").Value = 16
$ws.Range("AS30").Value = 21
$ws.Range("AT30").Value = 15
$ws.Range("AV30").Value = 20
$ws.Range("AX30").Value = 18
$ws.Range("BF30").Value = "'2008-11-20"

# Row 31
$ws.Range("AF31").Value = 25
$ws.Range("AJ31").Value = 16
$ws.Range("AU31").Value = 25
$ws.Range("AV31").Value = 13
$ws.Range("BA31").Value = 12
$ws.Range("BF31").Value = "'2008-11-20"
